# Applies the updated cryptocurrency price/volume snapshot to sheet1.
# Note: many "Price" values look like plain numbers (e.g. "18.70", "0.9992")
# but must stay EXACT text (Excel would otherwise silently drop trailing
# zeros / renormalize them as floats). We force text by prefixing the value
# with a leading apostrophe, PowerShell-escaped as two quotes (`''`) at the
# start of the single-quoted literal, e.g. '''18.70' => literal string `'18.70`.
# "Thousands-grouped" style prices (e.g. "30.405.53") are not valid numbers
# so Excel already keeps them as text without any prefix needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.405.53'
$ws.Range("E2").Value = '  -1.83%  '

# Row 3
$ws.Range("D3").Value = '1.905.24'
$ws.Range("E3").Value = '  -2.71%  '

# Row 4
$ws.Range("D4").Value = '''0.9992'
$ws.Range("E4").Value = '  -0.19%  '

# Row 5
$ws.Range("D5").Value = '''238.87'
$ws.Range("E5").Value = '  -2.09%  '

# Row 6
$ws.Range("D6").Value = '''0.9986'
$ws.Range("E6").Value = '  -0.26%  '

# Row 7
$ws.Range("D7").Value = '''0.4725'
$ws.Range("E7").Value = '  -2.44%  '

# Row 8
$ws.Range("D8").Value = '''0.2836'
$ws.Range("E8").Value = '  -3.56%  '

# Row 9
$ws.Range("D9").Value = '''0.06663'
$ws.Range("E9").Value = '  -6.15%  '

# Row 10
$ws.Range("D10").Value = '''18.70'
$ws.Range("E10").Value = '  -5.04%  '

# Row 11
$ws.Range("D11").Value = '''99.95'
$ws.Range("E11").Value = '  -6.56%  '

# Row 12
$ws.Range("D12").Value = '''0.07704'
$ws.Range("E12").Value = '  -0.72%  '

# Row 13
$ws.Range("D13").Value = '1.915.37'
$ws.Range("E13").Value = '  -2.34%  '

# Row 14
$ws.Range("D14").Value = '''5.199'
$ws.Range("E14").Value = '  -3.53%  '

# Row 15
$ws.Range("D15").Value = '''0.6663'
$ws.Range("E15").Value = '  -5.50%  '

# Row 16
$ws.Range("D16").Value = '30.410.02'
$ws.Range("E16").Value = '  -1.84%  '

# Row 17
$ws.Range("D17").Value = '''253.35'
$ws.Range("E17").Value = '  -9.09%  '

# Row 18
$ws.Range("D18").Value = '''0.9987'
$ws.Range("E18").Value = '  -0.24%  '

# Row 19
$ws.Range("D19").Value = '''0.000007440'
$ws.Range("E19").Value = '  -5.00%  '

# Row 20
$ws.Range("D20").Value = '''12.60'
$ws.Range("E20").Value = '  -5.39%  '

# Row 21
$ws.Range("D21").Value = '''5.369'
$ws.Range("E21").Value = '  -2.72%  '

# Row 22
$ws.Range("D22").Value = '''0.9991'
$ws.Range("E22").Value = '  -0.16%  '

# Row 23: Chainlink
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").Value = '''6.314'
$ws.Range("E23").Value = '  -3.11%  '

# Row 24: Cosmos
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '''9.391'
$ws.Range("E24").Value = '  -3.86%  '

# Row 25: Monero
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '''167.45'
$ws.Range("E25").Value = '  -1.11%  '

# Row 26: EthereumClassic
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = '''18.85'
$ws.Range("E26").Value = '  -4.47%  '

# Row 27: LidoDAOToken
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '''2.042'
$ws.Range("E27").Value = '  -6.29%  '

# Row 28: Stellar
$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D28").Value = '''0.1008'
$ws.Range("E28").Value = '  -4.13%  '

# Row 29: Filecoin
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").Value = '''4.657'
$ws.Range("E29").Value = '  +0.73%  '

# Row 30
$ws.Range("D30").Value = '''1.368'
$ws.Range("E30").Value = '  -2.85%  '

# Row 31: PancakeSwap
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '''1.509'
$ws.Range("E31").Value = '  -3.74%  '

# Row 32: InternetComputer(DFINITY)
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '''4.245'
$ws.Range("E32").Value = '  -3.79%  '

# Row 33: Hedera
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '''0.04710'
$ws.Range("E33").Value = '  -3.81%  '

# Row 34: ImmutableX
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '''0.7276'
$ws.Range("E34").Value = '  -3.45%  '

# Row 35: ARBITRUM
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '''1.108'
$ws.Range("E35").Value = '  -5.32%  '

# Row 36: Frax
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").Value = '''0.9979'
$ws.Range("E36").Value = '  -0.31%  '

# Row 37: HuobiToken
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '''2.697'
$ws.Range("E37").Value = '  -1.32%  '

# Row 38: VeChain
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.01911'
$ws.Range("E38").Value = '  -4.80%  '

# Row 39: MXToken
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '''2.592'
$ws.Range("E39").Value = '  -3.41%  '

# Row 40: FraxShare
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '''6.234'
$ws.Range("E40").Value = '  -4.31%  '

# Row 41: Aave
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '''73.58'
$ws.Range("E41").Value = '  -6.09%  '

# Row 42: RenderToken
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '''1.955'
$ws.Range("E42").Value = '  -8.06%  '

# Row 43: TrustWalletToken
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '''0.8586'
$ws.Range("E43").Value = '  -4.15%  '

# Row 44: Quant
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''105.43'
$ws.Range("E44").Value = '  -3.65%  '

# Row 45: PaxDollar
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '''0.9979'
$ws.Range("E45").Value = '  -0.30%  '

# Row 46: TheSandbox
$ws.Range("B46").Value = 'TheSandbox'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D46").Value = '''0.4220'
$ws.Range("E46").Value = '  -5.19%  '

# Row 47: Maker
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '''987.60'
$ws.Range("E47").Value = '  +0.28%  '

# Row 48
$ws.Range("D48").Value = '''7.372'
$ws.Range("E48").Value = '  -6.76%  '

# Row 49: Algorand
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '''0.1195'
$ws.Range("E49").Value = '  -4.34%  '

# Row 50: Elrond
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '''34.57'
$ws.Range("E50").Value = '  -3.92%  '

# Row 51: EnergySwap
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''8.754'
$ws.Range("E51").Value = '  -6.01%  '
